# reviews.xlsx: add a "Length" column (character count of the Review text)
# between "Review" (E) and "Deviation from the Mean" (F), pushing the
# deviation column to G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Figure out how many data rows exist before we touch anything.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count()

# Insert a fresh column at F — this shifts the existing "Deviation from the
# Mean" column (and its header style) from F to G automatically.
$ws.Columns.Item(6).Insert()

# New header, matching the bold/centered/bordered style used by the other
# header cells (Insert() already carried that style onto F1 for us, but set
# it explicitly too so the text goes in safely).
$ws.Range("F1").Value = "Length"

# Fill in the review length (character count) for every data row. Use
# ToCharArray().Count() rather than .Length so characters outside the BMP
# (emoji, etc., encoded as UTF-16 surrogate pairs) count once each instead
# of twice — matching the Python len() the reference data was built with.
for ($r = 2; $r -le $lastRow; $r++) {
    $review = $ws.Cells.Item($r, 5).Value()
    if ($review -eq $null) {
        $len = 0
    } else {
        $len = $review.ToCharArray().Count()
    }
    $ws.Cells.Item($r, 6).Value = $len
}
